# Replace the body rows (2..9) of Sheet1 with the new scraped citation data
# and extend the table through row 16 (all APSR v85 issue-1 articles except
# the book review and front/back matter).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = '“Erratum.” American Political Science Review 85, no. 1 (1991): 234–234. https://doi.org/10.2307/1962964.'
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""

# Row 3
$ws.Range("A3").Value = 'Abramson, Paul R., and Charles W. Ostrom. “Macropartisanship: An Empirical Reassessment.” American Political Science Review 85, no. 1 (1991): 181–92. https://doi.org/10.2307/1962884.'
$ws.Range("B3").Value = 'To evaluate the comparability of the Gallup and Michigan Survey Research Center measures for studying levels of partisanship among the U.S. electorate we compare the overtime distribution of partisanship and the correlates of partisanship using the results of Gallup surveys, the National Election Studies, and the General Social Surveys. Compared with the Gallup results, both the other two surveys reveal less short-term variation and also less total variation. Compared with the Gallup results, the National Election Studies partisanship results are less related to short-term electoral outcomes and do not appear to be strongly driven by short-term economic and political evaluations. Our analyses suggest that scholars should be cautious about using Gallup results to revise conclusions based upon analyses that employ the Michigan Survey Research Center party identification measure.'
$ws.Range("C3").Value = 'Michigan State University'
$ws.Range("D3").Value = 'Michigan State University'

# Row 4
$ws.Range("A4").Value = 'Rapoport, Ronald B., Walter J. Stone, and Alan I. Abramowitz. “Do Endorsements Matter? Group Influence in the 1984 Democratic Caucuses.” American Political Science Review 85, no. 1 (1991): 193–203. https://doi.org/10.2307/1962885.'
$ws.Range("B4").Value = 'Although candidates are very concerned about garnering prenomination group endorsements, there is virtually no research on the effects of such endorsements at the presidential level. In 1984, teachers'', women''s, and labor groups all endorsed Mondale for the Democratic nomination. Using a survey of Democratic caucus participants, we examine the effects of these endorsements on both candidate support and prenomination activity for the endorsed candidate. We find substantial effects of labor and teachers'' group endorsements on the candidate choice of their members, but virtually no effect of women''s groups. Labor union membership also had the greatest impact on campaign activity for Mondale.'
$ws.Range("C4").Value = 'College of William and Mary'
$ws.Range("D4").Value = 'University of Colorado / Emory University'

# Row 5
$ws.Range("A5").Value = 'Scholz, John T. “Cooperative Regulatory Enforcement and the Politics of Administrative Effectiveness.” American Political Science Review 85, no. 1 (1991): 115–36. https://doi.org/10.2307/1962881.'
$ws.Range("B5").Value = 'Even when political interests control bureaucratic outputs, the control of policy outcomes is complicated by trade-offs between controllable versus effective implementation strategies. I use a nested game framework to explain why a cooperative strategy can increase enforcement effectiveness in the narrow administrative game and why principal-agent control problems and collective action problems associated with the strategy lead policy beneficiaries to oppose the effective strategy in the broader political games. Analyses of state-level Occupational Safety and Health Administration enforcement provide evidence that cooperation does enhance the impact of enforcement in reducing workplace injury rates but that policy beneficiaries oppose and sabotage cooperation. The interactions between administrative effectiveness and interest group politics in this and other implementation situations require that both be analyzed simultaneously, and the nested game framework can provide a systematic approach to such analyses.'
$ws.Range("C5").Value = 'State University of New York Stony Brook'
$ws.Range("D5").Value = ""

# Row 6
$ws.Range("A6").Value = 'Cohen, Jeffrey E., Michael A. Krassa, and John A. Hamman. “The Impact of Presidential Campaigning on Midterm U.S. Senate Elections.” American Political Science Review 85, no. 1 (1991): 165–78. https://doi.org/10.2307/1962883.'
$ws.Range("B6").Value = 'The conventional wisdom about presidential campaigning in midterm Senate elections is that presidential efforts lack impact or have negative impact. We discuss conceptual problems with the conventional view and offer an alternative that views presidential campaigning as strategic. We test this alternative and find support for it. Further, we find that presidential campaign efforts have a positive impact on the vote through the mobilization of nonvoters. Finally, in a significant number of cases, presidential campaigning may have been the margin of victory for candidates of the president''s party. We discuss the implications of these findings on assessments of the president and relations with Congress during the second half of the term.'
$ws.Range("C6").Value = 'University of Illinois Urbana'
$ws.Range("D6").Value = 'University of Illinois Urbana / Southern Illinois University Carbondale'

# Row 7
$ws.Range("A7").Value = 'Mitchell, Timothy. “The Limits of the State: Beyond Statist Approaches and Their Critics.” American Political Science Review 85, no. 1 (1991): 77–96. https://doi.org/10.2307/1962879.'
$ws.Range("B7").Value = 'The state has always been difficult to define. Its boundary with
               society appears elusive, porous, and mobile. I argue that this elusiveness should not
               be overcome by sharper definitions, but explored as a clue to the state''s nature.
               Analysis of the literature shows that neither rejecting the state in favor of such
               concepts as the political system, nor “bringing it back in,” has dealt with this
               boundary problem. The former approach founders on it, the latter avoids it by a
               narrow idealism that construes the state-society distinction as an external relation
               between subjective and objective entities. A third approach, presented here, can
               account for both the salience of the state and its elusiveness. Reanalyzing evidence
               presented by recent theorists, state-society boundaries are shown to be distinctions
               erected internally, as an aspect of more complex power relations. Their appearance
               can be historically traced to technical innovations of the modern social order,
               whereby methods of organization and control internal to the social processes they
               govern create the effect of a state structure external to those
               processes.'
$ws.Range("C7").Value = 'New York University'
$ws.Range("D7").Value = ""

# Row 8
$ws.Range("A8").Value = 'Sinopoli, Richard C., and Nancy J. Hirschmann. “Feminism and Liberal Theory.” American Political Science Review 85, no. 1 (1991): 221–34. https://doi.org/10.2307/1962887.'
$ws.Range("B8").Value = 'In her article on “Freedom, Recognition, and Obligation: A Feminist Approach to Political Theory,” published in the December 1989 issue of this Review, Nancy J. Hirschmann argued that a feminist methodology could breathe new and useful life into liberal political theory, relieving it of its structural sexism. In this Controversy, Richard C. Sinopoli takes issue with key claims made by Hirschmann. In turn, Hirschmann elaborates her case.'
$ws.Range("C8").Value = 'University of California, Davis'
$ws.Range("D8").Value = 'Cornell University'

# Row 9
$ws.Range("A9").Value = 'Baron, David P. “A Spatial Bargaining Theory of Government Formation in Parliamentary Systems.” American Political Science Review 85, no. 1 (1991): 137–64. https://doi.org/10.2307/1962882.'
$ws.Range("B9").Value = 'A theory of government formation in parliamentary systems is developed from a model incorporating policy-oriented parties with spatial preferences and a formation process in which parties are selected in sequence to attempt to form a government. A government is formed when the policy proposed by the party selected is sustained on a vote of confidence. The equilibria identify the government and the policy it will implement and depend on the configuration of preferences and on the government formation process. For example, in a political system with two large parties and one small party and in which one of the two large parties will be selected to attempt to form a government, the government will be formed by a large party with the support of the small party. The policy will reflect the preferences of the small party but will be closer to the ideal point of the large party.'
$ws.Range("C9").Value = 'Stanford University'
$ws.Range("D9").Value = ""

# Row 10
$ws.Range("A10").Value = 'Epstein, Lee, and C. K. Rowland. “Debunking the Myth of Interest Group Invincibility in the Courts.” American Political Science Review 85, no. 1 (1991): 205–17. https://doi.org/10.2307/1962886.'
$ws.Range("B10").Value = 'Research on interest group litigation has provoked a reevaluation of the conventional wisdom about the study of pressure group activity and judicial politics. Nevertheless, the notion that interest groups are intrepid litigators that rarely lose to nongroup adversaries persists unchallenged and unscathed. We seek to determine if groups are, in fact, as invincible as the literature suggests. Several findings emerge that may undermine conventional wisdom about the relative efficacy of group-sponsored litigation. Most important is that groups are no more likely than nongroups to win, at least in U.S. District Courts. Based on this and other results, we draw a number of conclusions about interest group litigation and the direction into which future study might head.'
$ws.Range("C10").Value = 'Southern Methodist University'
$ws.Range("D10").Value = 'University of Kansas'

# Row 11
$ws.Range("A11").Value = 'Goldstein, Joshua S., and John R. Freeman. “U.S.-Soviet-Chinese Relations: Routine, Reciprocity, or Rational Expectations?” American Political Science Review 85, no. 1 (1991): 17–35. https://doi.org/10.2307/1962876.'
$ws.Range("B11").Value = 'International relations theorists disagree about whether great power behaviors reflect bureaucratic routine or reciprocity. Recently, some have suggested that these behaviors result from great powers'' rational expectations rather than from simple routine or reciprocity. The debate is flawed in several respects. The quasi-experimental studies of great power behavior suffer from specification and measurement errors. Furthermore, most studies of great power behavior focus exclusively on the superpowers, without adequately appreciating China''s role in world politics. We present an improved analysis that recognizes the potential effects of Chinese behavior and ameliorates the methodological flaws in existing work. The results indicate that the behaviors of the United States, the Soviet Union, and China are a relatively stable mix of bureaucratic routine and reciprocity. The results also indicate complex, asymmetrical connections among U.S.-Soviet, U.S.-Chinese, and Soviet-Chinese relations, consistent with the notion of a strategic triangle.'
$ws.Range("C11").Value = 'University of Southern California'
$ws.Range("D11").Value = 'University of Minnesota'

# Row 12
$ws.Range("A12").Value = 'Shklar, Judith N. “Redeeming American Political Theory.” American Political Science Review 85, no. 1 (1991): 3–15. https://doi.org/10.2307/1962875.'
$ws.Range("B12").Value = 'American political theory has been accused of being uniformly liberal; but its history is diverse and is worth studying to understand the development of political science and the institutions it reflects (representative government, federalism, judicial review, and slavery). While modern social science expresses a slow democratization of values, it has been compatible with many ideologies. This can be seen in Jefferson''s anthropology, Madison''s theory of collective rationality, and Hamilton''s empirical political economy. Jacksonian democracy encouraged social history, while its opponents devised an elitist political sociology. Southern defenders of slavery were the earliest to develop a deterministic and authoritarian sociology, but after the Civil War Northern thinkers emulated them with Social Darwinism and quests for causal laws to grasp constant change in industrial society. Though social critics abounded, democratic empirical theory emerged in the universities only in the generation of Merriam and Dewey, who founded contemporary political science.'
$ws.Range("C12").Value = 'Harvard University'
$ws.Range("D12").Value = ""

# Row 13
$ws.Range("A13").Value = 'Honig, B. “Declarations of Independence: Arendt and Derrida on the Problem of Founding a Republic.” American Political Science Review 85, no. 1 (1991): 97–113. https://doi.org/10.2307/1962880.'
$ws.Range("B13").Value = 'Beginning with Hannah Arendt''s depiction of the American Revolution and founding, I critically examine Arendt''s reading of the Declaration of Independence, comparing it with Jacques Derrida''s reading of (a draft of) the same document, in order to show that Arendt is careless in her easy dismissal of the declaration''s essentialist moments. Derrida, it seems to me, has a better, more subtle appreciation of the both necessary and impossible role of essentialism in modern political theory and practice. I conclude, however, that Arendt nonetheless succeeds in theorizing a powerful and suggestive practice of political authority for modernity, a practice that is uniquely activist and appropriate for a democratic politics.'
$ws.Range("C13").Value = 'Harvard University'
$ws.Range("D13").Value = ""

# Row 14
$ws.Range("A14").Value = 'Ostrom, Elinor. “Rational Choice Theory and Institutional Analysis: Toward Complementarity.” American Political Science Review 85, no. 1 (1991): 237–43. https://doi.org/10.2307/1962888.'
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = 'Indiana University'
$ws.Range("D14").Value = ""

# Row 15
$ws.Range("A15").Value = 'Schonhardt-Bailey, Cheryl. “Lessons in Lobbying for Free Trade in 19th-Century Britain: To Concentrate or Not.” American Political Science Review 85, no. 1 (1991): 37–58. https://doi.org/10.2307/1962877.'
$ws.Range("B15").Value = 'I present a modified version of the public choice interest group model that integrates concentrated and deconcentrated interests with successful lobbying. It is argued that effective free trade lobbying required the political fusion of the economic interests representing two fundamental changes in nineteenth-century Britain''s economy: (1) geographic concentration of the core export industry (cotton textiles) and (2) deconcentration of the broader export sector both geographically and in terms of industrial structure. Empirical evidence from both national and individual levels firmly supports the contention that the timing and political success of Britain''s nineteenth-century free trade lobby required the combined forces of core export interests and the more diverse and geographically more evenly distributed interests of the export sector as a whole.'
$ws.Range("C15").Value = 'University of California Los Angeles London School of Economics and Political Science'
$ws.Range("D15").Value = ""

# Row 16
$ws.Range("A16").Value = 'Booth, William James. “The New Household Economy.” American Political Science Review 85, no. 1 (1991): 59–75. https://doi.org/10.2307/1962878.'
$ws.Range("B16").Value = 'In the revival of scholarly interest in Marx''s relationship to one of the greatest of
his teachers, Aristotle, a crucial aspect has been undervalued: Marx''s indebtedness
to ancient political economy. While Marx employed the methods and key concepts of the
economic science of his day in analyzing capitalism, he embedded that explanation in
a higher-order theory of the economy. This theory, derived from the Aristotelian
account of the household economy, seeks to situate the economy in an overarching
account of the community, its purposiveness, and the place within it of activity,
time, and domination. Marx sought thereby to illuminate the historical uniqueness of
capitalism and, relatedly, to show the bounded character of the economic science
(including his own) appropriate to the understanding of it. Central elements of the
Aristotelian critique of an economy given over to Midas-like acquisition also find
their way into Marx''s evaluation of capitalism, and the ideal of the ancient
oikos forms one of the core parts of Marx''s theory of communism
as the new household economy.'
$ws.Range("C16").Value = 'McGill University'
$ws.Range("D16").Value = ""

